$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160
$r1.Borders.LineStyle = 1

$ws.Range("A2").Value = 0
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "disconnected_elements"
